$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: paragraph "As an authenticated user, I can see the aggr" +
# "egate results of my polls." (split across two runs by a _GoBack
# bookmark) becomes a single run with the complete sentence, and the
# bookmark is removed from here (it moves to the "delete polls" bullet
# below). A Find/Replace that spans the whole sentence (crossing the
# run/bookmark boundary) merges everything into one run automatically
# and drops the now-enclosed bookmark.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "As an authenticated user, I can see the aggregate results of my polls.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As an authenticated user, I can see the aggregate results of my polls.", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2 & 3: the "see and vote on everyone's polls" bullet and the
# "delete polls" bullet swap text content. The bullet that used to read
# "see and vote..." keeps its position but now reads "delete polls..."
# and picks up strikethrough formatting (plus the _GoBack bookmark at
# its end); the bullet that used to read "delete polls..." keeps its
# position too and now reads "see and vote...", unchanged formatting.
# ---------------------------------------------------------------------
$pVote = $d.Paragraphs.Item(6)

# sanity text captured defensively before either is modified
$voteText = "As an unauthenticated or authenticated user, I can see and vote on everyone's polls."
$deleteText = "As an authenticated user, I can delete polls that I decide I don't want anymore."

# --- paragraph 6: becomes the (struck-through) "delete polls" bullet ---
$r6 = $pVote.Range
$r6.MoveEnd(1, -1) | Out-Null
$r6.Text = $deleteText + "X"

$p6 = $d.Paragraphs.Item(6)
$full6 = $p6.Range
$full6.MoveEnd(1, -1) | Out-Null

# place the (still hidden/hyphenated) bookmark just before the trailing
# placeholder "X", then delete only the placeholder so the bookmark ends
# up sitting cleanly at the end of the run, right before the pilcrow.
$bmSpot = $d.Range($full6.End - 1, $full6.End - 1)
$d.Bookmarks.Add("_GoBack", $bmSpot) | Out-Null
$d.Range($full6.End - 1, $full6.End).Delete() | Out-Null

# apply strikethrough to both the run and the paragraph mark
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Font.StrikeThrough = 1

# --- paragraph 7: becomes the (plain) "see and vote" bullet ---
$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
$r7.MoveEnd(1, -1) | Out-Null
$r7.Text = $voteText
